# Update gh-pages to output generated at 456a3b4
# -------------------------------------------------------------
# The upstream scraper re-ran and refreshed the "想去人数" (interest
# count, column F) for a number of events across the 展览 / 演出 /
# 全部类型 sheets, and refreshed one event's cover image URL
# (column I) which had been reprocessed on bilibili's CDN.
# -------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ==== Sheet 展览 (exhibitions) ====
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 122
$ws1.Range("F3").Value  = 318
$ws1.Range("F5").Value  = 1188
$ws1.Range("F8").Value  = 134
$ws1.Range("F9").Value  = 137
$ws1.Range("F13").Value = 159
$ws1.Range("F14").Value = 1387
$ws1.Range("F15").Value = 525
$ws1.Range("F16").Value = 201
$ws1.Range("F17").Value = 314
$ws1.Range("F19").Value = 723
$ws1.Range("F20").Value = 1102
$ws1.Range("F23").Value = 2562
$ws1.Range("F24").Value = 1324
$ws1.Range("F26").Value = 233
$ws1.Range("F27").Value = 377
$ws1.Range("F28").Value = 983
$ws1.Range("F30").Value = 1124
$ws1.Range("F31").Value = 136
$ws1.Range("F33").Value = 760
$ws1.Range("F34").Value = 465
$ws1.Range("F35").Value = 601
$ws1.Range("F36").Value = 772
$ws1.Range("F37").Value = 330
$ws1.Range("F38").Value = 218
$ws1.Range("I35").Value = "//i2.hdslb.com/bfs/openplatform/202405/99kWb2dy1714964533903.png"

# ==== Sheet 演出 (performances) ====
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F13").Value = 546
$ws2.Range("F21").Value = 13

# ==== Sheet 全部类型 (all types) ====
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F6").Value  = 122
$ws4.Range("F7").Value  = 318
$ws4.Range("F11").Value = 1188
$ws4.Range("F14").Value = 134
$ws4.Range("F16").Value = 137
$ws4.Range("F20").Value = 159
$ws4.Range("F21").Value = 1387
$ws4.Range("F22").Value = 525
$ws4.Range("F23").Value = 201
$ws4.Range("F24").Value = 314
$ws4.Range("F26").Value = 1102
$ws4.Range("F27").Value = 2562
$ws4.Range("F29").Value = 1324
$ws4.Range("F34").Value = 233
$ws4.Range("F35").Value = 377
$ws4.Range("F36").Value = 983
$ws4.Range("F40").Value = 1124
$ws4.Range("F41").Value = 760
$ws4.Range("F42").Value = 465
$ws4.Range("F43").Value = 601
$ws4.Range("F44").Value = 772
$ws4.Range("F45").Value = 330
$ws4.Range("F47").Value = 13
$ws4.Range("F48").Value = 218
$ws4.Range("I43").Value = "//i2.hdslb.com/bfs/openplatform/202405/99kWb2dy1714964533903.png"
